$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Duplicate "Rainfall 2" (keeps its column widths, the A1 threaded
#    comment, number formats, etc.) to seed the new "Sheet1" tab that
#    will hold the B1 R3 bottle-weight data.
# ---------------------------------------------------------------------
$source = $wb.Worksheets.Item("Rainfall 2")
$null = $source.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$new = $wb.Worksheets.Item($wb.Worksheets.Count)
$new.Name = "Sheet1"

# ---------------------------------------------------------------------
# 2. The copied sheet still carries Rainfall-2's rain_date / collected
#    weight / formula / filter & extraction columns. Only
#    unique sample_id (A), sample_id (B), vol_water_applied (D),
#    bottle_mass (E) and filter_location (I) are populated for the new
#    B1 R3 batch, so blank out the rest of the data columns first.
# ---------------------------------------------------------------------
$new.Range("C2:C55").ClearContents()
$new.Range("E2:H55").ClearContents()
$new.Range("J2:M55").ClearContents()

# ---------------------------------------------------------------------
# 3. Write the new sample rows.
# ---------------------------------------------------------------------
$data = @(
    @(2, "SSCAMR00163", "C01R3", 88),
    @(3, "SSCAMR00164", "C02R3", 90),
    @(4, "SSCAMR00165", "C03R3", 89),
    @(5, "SSCAMR00166", "C04R3", 89),
    @(6, "SSCAMR00167", "C05R3", 89),
    @(7, "SSCAMR00168", "C06R3", 105),
    @(8, "SSCAMR00169", "C07R3", 90),
    @(9, "SSCAMR00170", "C08R3", 89),
    @(10, "SSCAMR00171", "C09R3", 88),
    @(11, "SSCAMR00172", "C10R3", 88),
    @(12, "SSCAMR00173", "C11R3", 87),
    @(13, "SSCAMR00174", "C12R3", 87),
    @(14, "SSCAMR00175", "C13R3", 88),
    @(15, "SSCAMR00176", "C14R3", 88),
    @(16, "SSCAMR00177", "C15R3", 89),
    @(17, "SSCAMR00178", "C16R3", 88),
    @(18, "SSCAMR00179", "C17R3", 88),
    @(19, "SSCAMR00180", "C18R3", 87),
    @(20, "SSCAMR00181", "C19R3", $null),
    @(21, "SSCAMR00182", "C20R3", $null),
    @(22, "SSCAMR00183", "C21R3", $null),
    @(23, "SSCAMR00184", "C22R3", $null),
    @(24, "SSCAMR00185", "C23R3", $null),
    @(25, "SSCAMR00186", "C24R3", $null),
    @(26, "SSCAMR00187", "C25R3", $null),
    @(27, "SSCAMR00188", "C26R3", $null),
    @(28, "SSCAMR00189", "C27R3", $null),
    @(29, "SSCAMR00190", "C28R3", $null),
    @(30, "SSCAMR00191", "C29R3", $null),
    @(31, "SSCAMR00192", "C30R3", $null),
    @(32, "SSCAMR00193", "C31R3", $null),
    @(33, "SSCAMR00194", "C32R3", $null),
    @(34, "SSCAMR00195", "C33R3", $null),
    @(35, "SSCAMR00196", "C34R3", $null),
    @(36, "SSCAMR00197", "C35R3", $null),
    @(37, "SSCAMR00198", "C36R3", $null),
    @(38, "SSCAMR00199", "C37R3", $null),
    @(39, "SSCAMR00200", "C38R3", $null),
    @(40, "SSCAMR00201", "C39R3", $null),
    @(41, "SSCAMR00202", "C40R3", $null),
    @(42, "SSCAMR00203", "C41R3", $null),
    @(43, "SSCAMR00204", "C42R3", $null),
    @(44, "SSCAMR00205", "C43R3", $null),
    @(45, "SSCAMR00206", "C44R3", $null),
    @(46, "SSCAMR00207", "C45R3", $null),
    @(47, "SSCAMR00208", "C46R3", $null),
    @(48, "SSCAMR00209", "C47R3", $null),
    @(49, "SSCAMR00210", "C48R3", $null),
    @(50, "SSCAMR00211", "C49R3", $null),
    @(51, "SSCAMR00212", "C50R3", $null),
    @(52, "SSCAMR00213", "C51R3", $null),
    @(53, "SSCAMR00214", "C52R3", $null),
    @(54, "SSCAMR00215", "C53R3", $null),
    @(55, "SSCAMR00216", "C54R3", $null),
)

foreach ($row in $data) {
    $new.Cells.Item($row[0], 1).Value = $row[1]
}
foreach ($row in $data) {
    $new.Cells.Item($row[0], 2).Value = $row[2]
}
foreach ($row in $data) {
    $new.Cells.Item($row[0], 4).Value = 1000
}
foreach ($row in $data) {
    $bottleMass = $row[3]
    if ($null -ne $bottleMass) {
        $new.Cells.Item($row[0], 5).Value = $bottleMass
    }
}
foreach ($row in $data) {
    $new.Cells.Item($row[0], 9).Value = "Freezer_F_Shelf_5"
}

# vol_water_applied(mL) column uses a thousands-separator number format.
$new.Range("D2:D55").NumberFormat = "#,##0"

# ---------------------------------------------------------------------
# 4. Selection / active-tab bookkeeping to match the saved workbook view
#    (the new sheet becomes the active tab, parked at E20; Rainfall 2
#    loses its old "in progress" cursor position).
# ---------------------------------------------------------------------
$null = $source.Select()
$null = $source.Cells.Select()
$null = $new.Select()
$null = $new.Range("E20").Select()
